# Painting Scavenger Hunt "system parameters" workbook edit
#
# Splits the former combined "hunt and painting" / "start-up" XML spec
# entries into separate, more granular rows:
#   - Specification Files: start_up -> menu, about, hunt, painting (4 rows
#     replacing the old single "start_up" + "hunt" rows, objects_of_interest
#     and end_goal shift down to make room)
#   - Hunt Parameters: source file renamed from "hunt and painting.xml" to
#     "hunt.xml"
#   - Menu Parameters: promoted to its own bold header row (like the other
#     sections) and its source file renamed from "start-up.xml" to "menu.xml"
#
# The sheet is laid out as a sequence of "sections", each starting with a
# bold/italic 14pt header row, and separated from the next section by a
# blank spacer row. We grow each inter-section gap from 1 row to 3 rows
# (net +2 each) and grow the first section's body from 4 rows to 6 rows
# (net +2), which is what pushes every later section down by the right
# amount (the same +2 increments stack: +2, +4, +4, +4 -> final offsets of
# +4, +6, +8 for Asset Directories / Hunt Parameters / Menu Parameters
# respectively).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the extra rows needed, bottom-to-top so the row numbers we
#        reference below (taken from the ORIGINAL layout) stay valid for
#        every subsequent Insert call. Each gap grows by 2 rows, and the
#        first section's body grows by 2 rows. ---

# Grow the gap before "Menu Parameters" (old row 19) from 1 -> 3 blank rows.
$ws.Rows.Item(18).Resize(2).Insert()

# Grow the gap before "Hunt Paramters" (old row 14) from 1 -> 3 blank rows.
$ws.Rows.Item(13).Resize(2).Insert()

# Grow the gap before "Asset Directories" (old row 7) from 1 -> 3 blank rows.
$ws.Rows.Item(6).Resize(2).Insert()

# Grow the "Specification Files" body (old rows 2-5) from 4 -> 6 rows.
$ws.Rows.Item(2).Resize(2).Insert()

# --- 2. The Insert() calls above copied the formatting of the row above
#        into the newly-created blank rows, which leaves stray styled-but-
#        empty cells sitting in what should be untouched gap rows. Clear
#        them completely so those rows disappear from the saved XML again
#        (matching the source, where gap rows have no <row> element). ---
$ws.Range("A8:B9").Clear()
$ws.Range("A17:B18").Clear()
$ws.Range("A24:B25").Clear()

# --- 3. Re-home the "Specification Files" section body (now rows 2-7) with
#        the correct labels/values. Force the normal (non-header) styling
#        explicitly, since a couple of these rows inherited the bold-italic
#        14pt header look from the Insert() above. ---
$ws.Range("A2:A7").Font.Bold = $true
$ws.Range("A2:A7").Font.Italic = $false
$ws.Range("A2:A7").Font.Size = 11

$ws.Range("B2:B7").Font.Bold = $false
$ws.Range("B2:B7").Font.Italic = $false
$ws.Range("B2:B7").Font.Size = 11
$ws.Range("B2:B7").NumberFormat = "@"

$ws.Range("A2").Value = "menu"
$ws.Range("B2").Value = "xml/menu.xml"

$ws.Range("A3").Value = "about"
$ws.Range("B3").Value = "xml/about.xml"

$ws.Range("A4").Value = "hunt"
$ws.Range("B4").Value = "xml/hunt.xml"

$ws.Range("A5").Value = "painting"
$ws.Range("B5").Value = "xml/painting.xml"

$ws.Range("A6").Value = "objects_of_interest"
$ws.Range("B6").Value = "xml/objects of interest.xml"

$ws.Range("A7").Value = "end_goal"
$ws.Range("B7").Value = "xml/end goal.xml"

# --- 4. "Hunt Paramters" section (now starting at row 20): the source
#        file moves from "hunt and painting.xml" to "hunt.xml". ---
$ws.Range("B20").Value = "hunt.xml"

# --- 5. "Menu Parameters" section (now starting at row 27): promote the
#        label to its own bold 14pt (non-italic) header row, matching the
#        look of the other section headers, and rename its source file
#        from "start-up.xml" to "menu.xml". ---
$ws.Rows.Item(27).RowHeight = 18.75
$ws.Range("A27").Value = "Menu Parameters"
$ws.Range("A27").Font.Bold = $true
$ws.Range("A27").Font.Italic = $false
$ws.Range("A27").Font.Size = 14

$ws.Range("B27").Value = "menu.xml"

# --- 6. Restore the active selection to where the author left off. ---
[void]$ws.Range("C9").Select()
